$wb = $excel.ActiveWorkbook

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1027.9767
$ws.Range("I129").Value = 310
$ws.Range("J129").Value = 1122.4474
$ws.Range("K129").Value = 930
$ws.Range("L129").Value = 3367.3422
$ws.Range("M129").Value = 4070
$ws.Range("N129").Value = -13367.3422

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4657.6895
$ws.Range("I137").Value = 4849.55
$ws.Range("J137").Value = 4231.3335
$ws.Range("K137").Value = 14548.65
$ws.Range("L137").Value = 12694.0005
$ws.Range("M137").Value = -11998.65
$ws.Range("N137").Value = -17794.0005

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 621
$ws.Range("I2").Value = 562.7692
$ws.Range("K2").Value = 562.7692
$ws.Range("M2").Value = -449.7692

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6524.3794
$ws.Range("I32").Value = 4611.0527
$ws.Range("K32").Value = 4611.0527
$ws.Range("M32").Value = -4324.0527

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3778.5642
$ws.Range("I74").Value = 3669.2727
$ws.Range("J74").Value = 4379.6665
$ws.Range("K74").Value = 3669.2727
$ws.Range("L74").Value = 4379.6665
$ws.Range("M74").Value = -2795.2727
$ws.Range("N74").Value = -6127.6665

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3778.5642
$ws.Range("I77").Value = 3669.2727
$ws.Range("J77").Value = 4379.6665
$ws.Range("K77").Value = 18346.3635
$ws.Range("L77").Value = 21898.3325
$ws.Range("M77").Value = -13978.3635
$ws.Range("N77").Value = -30634.3325

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 621
$ws.Range("I116").Value = 562.7692
$ws.Range("K116").Value = 562.7692
$ws.Range("M116").Value = 1731.2308

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3589
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 4355.7144
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 13067.1432
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -17967.1432

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 621
$ws.Range("I3").Value = 562.7692
$ws.Range("K3").Value = 562.7692
$ws.Range("M3").Value = -448.7692

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1390.9375
$ws.Range("I107").Value = 1301.3
$ws.Range("J107").Value = 1540.3334
$ws.Range("K107").Value = 1301.3
$ws.Range("L107").Value = 1540.3334
$ws.Range("M107").Value = 618.7
$ws.Range("N107").Value = -5380.3334

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 16671624
$ws.Range("I99").Value = 40003220
$ws.Range("J99").Value = 6198.5713
$ws.Range("K99").Value = 40003220
$ws.Range("L99").Value = 6198.5713
$ws.Range("M99").Value = -40001722
$ws.Range("N99").Value = -9194.5713

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 16671624
$ws.Range("I126").Value = 40003220
$ws.Range("J126").Value = 6198.5713
$ws.Range("K126").Value = 120009660
$ws.Range("L126").Value = 18595.7139
$ws.Range("M126").Value = -120007190
$ws.Range("N126").Value = -23535.7139

# Sheet CUL, row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3726.1875
$ws.Range("I3").Value = 2514.9167
$ws.Range("J3").Value = 7360
$ws.Range("K3").Value = 7544.750100000001
$ws.Range("L3").Value = 22080
$ws.Range("M3").Value = -7432.750100000001
$ws.Range("N3").Value = -22304

# Sheet CUL, row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 7200.6665
$ws.Range("I92").Value = 800
$ws.Range("J92").Value = 10401
$ws.Range("K92").Value = 2400
$ws.Range("L92").Value = 31203
$ws.Range("M92").Value = -1152
$ws.Range("N92").Value = -33699

# Sheet CUL, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 549.1818
$ws.Range("J107").Value = 764.2222
$ws.Range("L107").Value = 2292.6666
$ws.Range("N107").Value = -6132.6666

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 695.35297
$ws.Range("J113").Value = 695.8182
$ws.Range("L113").Value = 2087.4546
$ws.Range("N113").Value = -6427.4546

# Sheet CUL, row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3887.2222
$ws.Range("I138").Value = 5750
$ws.Range("J138").Value = 3654.375
$ws.Range("K138").Value = 17250
$ws.Range("L138").Value = 10963.125
$ws.Range("M138").Value = -12110
$ws.Range("N138").Value = -21243.125

# Sheet CUL, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2385
$ws.Range("I139").Value = 1423.75
$ws.Range("J139").Value = 3666.6667
$ws.Range("K139").Value = 4271.25
$ws.Range("L139").Value = 11000.0001
$ws.Range("M139").Value = 868.75
$ws.Range("N139").Value = -21280.0001

# Sheet GSM, row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13334.389
$ws.Range("I43").Value = 1669.0834
$ws.Range("J43").Value = 36665
$ws.Range("K43").Value = 1669.0834
$ws.Range("L43").Value = 36665
$ws.Range("M43").Value = -1518.0834
$ws.Range("N43").Value = -36967

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2367
$ws.Range("I102").Value = 1480.0834
$ws.Range("J102").Value = 4495.6
$ws.Range("K102").Value = 1480.0834
$ws.Range("L102").Value = 4495.6
$ws.Range("M102").Value = 141.9166
$ws.Range("N102").Value = -7739.6

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10321.571
$ws.Range("I40").Value = 10900.5
$ws.Range("K40").Value = 10900.5
$ws.Range("M40").Value = -10764.5

# Sheet LTW, row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 304.6842
$ws.Range("I55").Value = 208.91667
$ws.Range("J55").Value = 468.85715
$ws.Range("K55").Value = 208.91667
$ws.Range("L55").Value = 468.85715
$ws.Range("M55").Value = -35.91667000000001
$ws.Range("N55").Value = -814.85715

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7109.3335
$ws.Range("I122").Value = 4062.4
$ws.Range("J122").Value = 9285.714
$ws.Range("K122").Value = 12187.2
$ws.Range("L122").Value = 27857.142
$ws.Range("M122").Value = -9737.200000000001
$ws.Range("N122").Value = -32757.142

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3904.3333
$ws.Range("I132").Value = 1934.375
$ws.Range("J132").Value = 5274.7393
$ws.Range("K132").Value = 5803.125
$ws.Range("L132").Value = 15824.2179
$ws.Range("M132").Value = -3273.125
$ws.Range("N132").Value = -20884.2179

# Sheet WVR, row 42
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 43365.668
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 43365.668
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 43365.668
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -44121.668

# Sheet WVR, row 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 18007.8
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 28013
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 28013
$ws.Range("M43").Value = -2851
$ws.Range("N43").Value = -28311

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5953695.5
$ws.Range("I81").Value = 8929919
$ws.Range("J81").Value = 1249.75
$ws.Range("K81").Value = 17859838
$ws.Range("L81").Value = 2499.5
$ws.Range("M81").Value = -17858777
$ws.Range("N81").Value = -4621.5

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5953695.5
$ws.Range("I84").Value = 8929919
$ws.Range("J84").Value = 1249.75
$ws.Range("K84").Value = 89299190
$ws.Range("L84").Value = 12497.5
$ws.Range("M84").Value = -89293886
$ws.Range("N84").Value = -23105.5
